# Partial commit of second revision
# Adds two new ObjectLocator rows ("btnEdit" / "txtRprAmount") to the
# RecurringPayment sheet, matching the style of the surrounding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RecurringPayment")
$ws.Activate() | Out-Null

# --- Row 28: btnEdit / by_id / btnEdit -------------------------------------
# Copy formatting from row 21 (same "code" style pattern: A/C use style 9).
$ws.Range("A21:C21").Copy($ws.Range("A28:C28")) | Out-Null
$ws.Range("A28").Value = "btnEdit"
$ws.Range("B28").Value = "by_id"
$ws.Range("C28").Value = "btnEdit"

# --- Row 29: txtRprAmount / by_xpath / xpath expression ---------------------
# Copy formatting from row 23 (plain, unstyled row).
# (The xpath text is written before the txtRprAmount label so the new
# shared-string entries land in the same order as the authored workbook.)
$ws.Range("A23:C23").Copy($ws.Range("A29:C29")) | Out-Null
$ws.Range("B29").Value = "by_xpath"
$ws.Range("C29").Value = "//input[contains(@id,'vamount')]/preceding-sibling::input"
$ws.Range("A29").Value = "txtRprAmount"

# Column C needs to widen to fit the new, longer xpath text.
$ws.Columns.Item(3).ColumnWidth = 65.5

# Leave the selection where the author's session left it.
$ws.Range("A36").Select() | Out-Null

Write-Output "done"
